$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the previous data row (row 10) down into the new row 11
$ws.Range("A10:G10").Copy() | Out-Null
$ws.Range("A11:G11").PasteSpecial(-4122) | Out-Null

# Fill in the new pin-mapping row (#10 in the "No." column)
$ws.Range("A11").Value = 10
$ws.Range("B11").Value = "PA8"
$ws.Range("C11").Value = "DO"
$ws.Range("D11").Value = "结果输出"
$ws.Range("E11").Value = "O"
$ws.Range("F11").ClearContents() | Out-Null
$ws.Range("G11").ClearContents() | Out-Null

# Match the author's last selected cell in the sheet view
$ws.Range("D19").Select() | Out-Null
